# Generate Report for Archive
#
# The 591498c8-...md and c82cffd3-...md entries swap their row positions
# (rows 4 and 5) on every sheet, and the 435f0468-...md entry (row 3) moves
# from "Ready for handoff" to "In Translation" on every sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Row 3 (435f0468...): status text only changes.
$ws.Range("E3").Value = "In Translation"
$ws.Range("F3").Value = "In Translation"

# Row 4 currently holds 591498c8..., row 5 currently holds c82cffd3... .
# Swap the two rows' data completely (File Name / Path And Name / Status x2 / Date).
$ws.Range("A4").Value = "c82cffd3-b3c4-43ac-8860-d58d8f741137.md"
$ws.Range("B4").Value = "e2e\c82cffd3-b3c4-43ac-8860-d58d8f741137.md"
$ws.Range("E4").Value = "In Translation"
$ws.Range("F4").Value = "In Translation"
$ws.Range("G4").Value = "2016-08-31 02:45:36"

$ws.Range("A5").Value = "591498c8-17e4-4d9a-ae5e-58e1bdbeed2f.md"
$ws.Range("B5").Value = "e2e\591498c8-17e4-4d9a-ae5e-58e1bdbeed2f.md"
$ws.Range("E5").Value = "Ready for handoff"
$ws.Range("F5").Value = "Ready for handoff"
$ws.Range("G5").Value = "2016-08-31 02:44:19"

# Hyperlinks stay anchored to their cell (and keep pointing at the same
# external target), only their displayed text follows the new cell value.
foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$B$4') {
        $h.TextToDisplay = "e2e\c82cffd3-b3c4-43ac-8860-d58d8f741137.md"
    }
    elseif ($addr -eq '$B$5') {
        $h.TextToDisplay = "e2e\591498c8-17e4-4d9a-ae5e-58e1bdbeed2f.md"
    }
}

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

# Row 3 (435f0468...): status text only changes.
$ws.Range("C3").Value = "In Translation"

# Row 4 currently holds 591498c8..., row 5 currently holds c82cffd3... .
$ws.Range("A4").Value = "c82cffd3-b3c4-43ac-8860-d58d8f741137.md"
$ws.Range("G4").Value = "c82cffd3-b3c4-43ac-8860-d58d8f741137.e608ac403c581f6fe0634a78dde703f7e6013893.zh-cn.xlf"
$ws.Range("H4").Value = "2016-08-31 02:45:31"

$ws.Range("A5").Value = "591498c8-17e4-4d9a-ae5e-58e1bdbeed2f.md"
$ws.Range("G5").Value = "591498c8-17e4-4d9a-ae5e-58e1bdbeed2f.9d8ec1e3e1d2c93bb33f2632da63a3790228b47b.zh-cn.xlf"
$ws.Range("H5").Value = "2016-08-31 02:44:14"

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$4') {
        $h.TextToDisplay = "c82cffd3-b3c4-43ac-8860-d58d8f741137.md"
    }
    elseif ($addr -eq '$A$5') {
        $h.TextToDisplay = "591498c8-17e4-4d9a-ae5e-58e1bdbeed2f.md"
    }
}

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

# Row 3 (435f0468...): status text only changes.
$ws.Range("C3").Value = "In Translation"

# Row 4 currently holds 591498c8..., row 5 currently holds c82cffd3... .
$ws.Range("A4").Value = "c82cffd3-b3c4-43ac-8860-d58d8f741137.md"
$ws.Range("G4").Value = "c82cffd3-b3c4-43ac-8860-d58d8f741137.e608ac403c581f6fe0634a78dde703f7e6013893.de-de.xlf"
$ws.Range("H4").Value = "2016-08-31 02:45:36"

$ws.Range("A5").Value = "591498c8-17e4-4d9a-ae5e-58e1bdbeed2f.md"
$ws.Range("G5").Value = "591498c8-17e4-4d9a-ae5e-58e1bdbeed2f.9d8ec1e3e1d2c93bb33f2632da63a3790228b47b.de-de.xlf"
$ws.Range("H5").Value = "2016-08-31 02:44:19"

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$4') {
        $h.TextToDisplay = "c82cffd3-b3c4-43ac-8860-d58d8f741137.md"
    }
    elseif ($addr -eq '$A$5') {
        $h.TextToDisplay = "591498c8-17e4-4d9a-ae5e-58e1bdbeed2f.md"
    }
}
